$wb = $excel.ActiveWorkbook

# Excel constant values
$xlEdgeLeft      = 7
$xlEdgeRight     = 10
$xlContinuous    = 1
$xlNone          = -4142
$xlPasteFormats  = -4122

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Build two helper/template cells (off to the side, will be cleared afterwards)
# that carry the two new border styles we need:
#   helperA -> top+bottom border only (no left/right)   -> used for the "C" columns
#   helperB -> top+bottom+right border (no left)        -> used for the "D" columns
$helperA = $ws2.Range("Z1")
$helperA.Style = "Normal"
$helperA.Borders.LineStyle = $xlContinuous
$helperA.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
$helperA.Borders.Item($xlEdgeRight).LineStyle = $xlNone

$helperB = $ws2.Range("Z2")
$helperB.Style = "Normal"
$helperB.Borders.LineStyle = $xlContinuous
$helperB.Borders.Item($xlEdgeLeft).LineStyle = $xlNone

# Apply the template formats (single paste operation per target cell keeps the
# workbook's style table minimal/deterministic instead of rebuilding borders
# edge-by-edge on every target cell)
$helperA.Copy()
$ws1.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$helperB.Copy()
$ws1.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

# Remove the helper cells, they must not remain part of the data
$helperA.Clear()
$helperB.Clear()

# --- Rename "fedcore" -> "approach" in the header rows of both sheets
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty cell G5 on the computational_comparison sheet
$ws2.Range("G5").ClearContents()
